$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.019.08"

$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("D5").Value = "'2.13"
$ws.Range("E5").Value = "  +12.35%  "

$ws.Range("D6").Value = "'235.28"

$ws.Range("D7").Value = "'654.82"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").Value = "'0.437"
$ws.Range("E8").Value = "  +2.23%  "

$ws.Range("D9").Value = "'1.10"
$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("D10").Value = "'0.999"
$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("D11").Value = "3.684.68"
$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("D12").Value = "'0.0000310"
$ws.Range("E12").Value = "  +14.91%  "

$ws.Range("D13").Value = "'44.42"
$ws.Range("E13").Value = "  -2.46%  "

$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("E15").Value = "  -1.31%  "

$ws.Range("D16").Value = "4.375.41"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").Value = "96.795.25"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "'8.57"
$ws.Range("E18").Value = "  -5.85%  "

$ws.Range("D19").Value = "3.691.40"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").Value = "'18.51"
$ws.Range("E21").Value = "  -4.57%  "

$ws.Range("E22").Value = "  +1.46%  "

$ws.Range("D23").Value = "'517.22"
$ws.Range("E23").Value = "  -1.64%  "

$ws.Range("E24").Value = "  -2.92%  "

$ws.Range("D25").Value = "'0.0000220"

$ws.Range("E26").Value = "  -3.86%  "

$ws.Range("D27").Value = "'110.75"
$ws.Range("E27").Value = "  +8.62%  "

$ws.Range("D28").Value = "'0.202"
$ws.Range("E28").Value = "  +19.74%  "

$ws.Range("E29").Value = "  -0.76%  "

$ws.Range("D30").Value = "'12.53"
$ws.Range("E30").Value = "  -1.23%  "

$ws.Range("E31").Value = "  -2.10%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  -4.09%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("D36").Value = "'32.54"
$ws.Range("E36").Value = "  -0.98%  "

$ws.Range("E37").Value = "  -3.24%  "

$ws.Range("D38").Value = "'629.23"
$ws.Range("E38").Value = "  -4.88%  "

$ws.Range("D39").Value = "'8.70"
$ws.Range("E39").Value = "  -3.93%  "

$ws.Range("E41").Value = "  +1.54%  "

$ws.Range("D42").Value = "'6.80"
$ws.Range("E42").Value = "  -5.76%  "

$ws.Range("D43").Value = "'0.492"
$ws.Range("E43").Value = "  +9.15%  "

$ws.Range("D44").Value = "'2.00"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").Value = "'39.99"
$ws.Range("E45").Value = "  +1.79%  "

$ws.Range("D46").Value = "'0.950"
$ws.Range("E46").Value = "  -2.43%  "

$ws.Range("E47").Value = "  -3.46%  "

$ws.Range("E48").Value = "  +1.47%  "

$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("E50").Value = "  -1.89%  "

$ws.Range("D51").Value = "'3.31"
$ws.Range("E51").Value = "  +2.43%  "

